# Rewrites the daily-job summary table (A2:E42) to the target values.
# The underlying data changed (a placeholder "ColumnN"/"Column1" row was
# inserted per day, plus two new rows for 2025-09-13), so every data row
# is written explicitly here rather than patched in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A keeps the workbook's existing date format (YYYY-MM-DD, same
# style used by the current rows) for both existing and newly-added rows.
$dateFmt = $ws.Cells.Item(2, 1).NumberFormat

$rows = @(
    @{ r = 2; date = 45908; job = 'ARA3A'; cnt = 5; hrs = 47.5; emp = 'Benito A, Daniel GS, Evaristo A, Leobardo RL, Omar S' },
    @{ r = 3; date = 45908; job = 'Column8'; cnt = 1; hrs = 4; emp = 'Column1' },
    @{ r = 4; date = 45908; job = 'Founders 2'; cnt = 3; hrs = 28.5; emp = 'Honorio G, Rigoberto Al-B, Rogelio M' },
    @{ r = 5; date = 45908; job = 'GU Henle'; cnt = 4; hrs = 38; emp = 'Fernando V, Laurentino, Noe VL, Oscar VS' },
    @{ r = 6; date = 45908; job = 'HanoverSpring'; cnt = 6; hrs = 57; emp = 'Antoine F, Carlos Al-V, Gaudencio B, Jose P, Moises P, Rata F' },
    @{ r = 7; date = 45908; job = 'Rowan'; cnt = 8; hrs = 76; emp = 'Adalberto T, Daniel LG, Elvis T, Henry G, Julio M, Luis Enrique R, Luis Martin R, Trinidad T' },
    @{ r = 8; date = 45908; job = 'Tidal Basin'; cnt = 6; hrs = 57; emp = 'Alberto R, Danis BA, Eduardo H, Feliciano R, Isidro M, Juan HR' },
    @{ r = 9; date = 45908; job = 'Wardman'; cnt = 14; hrs = 133; emp = 'Alfonso D, Andres G, Benny S, Carlos G, Cristobal L, Diego R, Eliacim R, Jesus L, Jose Luis H, Juan G, Miguel A, Misael M, Pablo G, William A' },
    @{ r = 10; date = 45909; job = 'ARA3A     Moorefield'; cnt = 7; hrs = 66.5; emp = 'Benito A, Daniel GS, Eric M R, Evaristo A, Henry G, Leobardo RL, Omar S' },
    @{ r = 11; date = 45909; job = 'Canvas'; cnt = 1; hrs = 9.5; emp = 'Noe VL' },
    @{ r = 12; date = 45909; job = 'Column13'; cnt = 1; hrs = 9; emp = 'Column1' },
    @{ r = 13; date = 45909; job = 'Founders 2'; cnt = 3; hrs = 27; emp = 'Honorio G, Rigoberto Al-B, Rogelio M' },
    @{ r = 14; date = 45909; job = 'GU Henle'; cnt = 3; hrs = 28.5; emp = 'Fernando V, Laurentino, Oscar VS' },
    @{ r = 15; date = 45909; job = 'HanoverSpring'; cnt = 5; hrs = 47.5; emp = 'Carlos Al-V, Gaudencio B, Jose P, Moises P, Rata F' },
    @{ r = 16; date = 45909; job = 'Rowan'; cnt = 7; hrs = 66.5; emp = 'Adalberto T, Daniel LG, Elvis T, Julio M, Luis Enrique R, Luis Martin R, Trinidad T' },
    @{ r = 17; date = 45909; job = 'Tidal Basin'; cnt = 6; hrs = 57; emp = 'Alberto R, Danis BA, Eduardo H, Feliciano R, Isidro M, Juan HR' },
    @{ r = 18; date = 45909; job = 'Wardman'; cnt = 12; hrs = 114; emp = 'Alfonso D, Andres G, Benny S, Carlos G, Cristobal L, Diego R, Eliacim R, Jesus L, Juan G, Miguel A, Pablo G, William A' },
    @{ r = 19; date = 45910; job = '2011 Crystal'; cnt = 3; hrs = 28.5; emp = 'Alejandro M S, Rigoberto Al-B, Rogelio M' },
    @{ r = 20; date = 45910; job = 'BridgeDist'; cnt = 2; hrs = 19; emp = 'Elvis T, Evaristo A' },
    @{ r = 21; date = 45910; job = 'Column18'; cnt = 1; hrs = 14; emp = 'Column1' },
    @{ r = 22; date = 45910; job = 'GU Henle'; cnt = 5; hrs = 49.5; emp = 'Fernando V, Jose Carlos G, Laurentino, Noe VL, Oscar VS' },
    @{ r = 23; date = 45910; job = 'HanoverSpring'; cnt = 4; hrs = 38; emp = 'Carlos Al-V, Jose P, Moises P, Rata F' },
    @{ r = 24; date = 45910; job = 'Kingstowne'; cnt = 5; hrs = 48; emp = 'Adalberto T, Gaudencio B, Luis Enrique R, Luis Martin R, Trinidad T' },
    @{ r = 25; date = 45910; job = 'Moorfield'; cnt = 7; hrs = 67; emp = 'Benito A, Daniel GS, Daniel LG, Eric M R, Julio M, Leobardo RL, Omar S' },
    @{ r = 26; date = 45910; job = 'Tidal Basin'; cnt = 6; hrs = 48; emp = 'Alberto R, Danis BA, Eduardo H, Feliciano R, Isidro M, Juan HR' },
    @{ r = 27; date = 45910; job = 'Wardman'; cnt = 12; hrs = 114; emp = 'Alfonso D, Andres G, Benny S, Carlos G, Cristobal L, Diego R, Eliacim R, Jesus L, Juan G, Miguel A, Pablo G, William A' },
    @{ r = 28; date = 45911; job = '2011 Crystal'; cnt = 3; hrs = 27; emp = 'Alejandro M S, Gerardo D, Rigoberto Al-B' },
    @{ r = 29; date = 45911; job = '2011 Crystal    Yard'; cnt = 1; hrs = 9; emp = 'Rogelio M' },
    @{ r = 30; date = 45911; job = 'Column23'; cnt = 1; hrs = 19; emp = 'Column1' },
    @{ r = 31; date = 45911; job = 'GU Henle'; cnt = 5; hrs = 47.5; emp = 'Fernando V, Jose Carlos G, Laurentino, Noe VL, Oscar VS' },
    @{ r = 32; date = 45911; job = 'HanoverSpring'; cnt = 5; hrs = 47.5; emp = 'Antoine F, Carlos Al-V, Jose P, Moises P, Rata F' },
    @{ r = 33; date = 45911; job = 'Kingstowne'; cnt = 8; hrs = 76; emp = 'Adalberto T, Elvis T, Gaudencio B, Julio M, Leobardo RL, Luis Enrique R, Luis Martin R, Trinidad T' },
    @{ r = 34; date = 45911; job = 'Moorefield'; cnt = 6; hrs = 57; emp = 'Benito A, Daniel GS, Daniel LG, Eric M R, Evaristo A, Omar S' },
    @{ r = 35; date = 45911; job = 'Tidal Basin'; cnt = 7; hrs = 66.5; emp = 'Alberto R, Danis BA, Eduardo H, Feliciano R, Isidro M, Juan HR, William A' },
    @{ r = 36; date = 45911; job = 'Wardman'; cnt = 12; hrs = 111; emp = 'Alfonso D, Andres G, Benny S, Carlos G, Cristobal L, Diego R, Eliacim R, Jesus L, Juan G, Miguel A, Misael M, Pablo G' },
    @{ r = 37; date = 45912; job = 'Column28'; cnt = 1; hrs = 24; emp = 'Column1' },
    @{ r = 38; date = 45912; job = 'Tidal Basin'; cnt = 1; hrs = 9.5; emp = 'Eduardo H' },
    @{ r = 39; date = 45912; job = 'Wardman'; cnt = 9; hrs = 82; emp = 'Alfonso D, Andres G, Cristobal L, Eliacim R, Jesus L, Juan G, Miguel A, Misael M, Pablo G' },
    @{ r = 40; date = 45913; job = 'Column28'; cnt = 1; hrs = 29; emp = 'Column1' },
    @{ r = 41; date = 45913; job = 'Tidal Basin'; cnt = 1; hrs = 4; emp = 'Eduardo H' },
    @{ r = 42; date = 45913; job = 'Wardman'; cnt = 9; hrs = 49.5; emp = 'Alfonso D, Andres G, Cristobal L, Eliacim R, Jesus L, Juan G, Miguel A, Misael M, Pablo G' }
)

foreach ($row in $rows) {
    $cellA = $ws.Cells.Item($row.r, 1)
    $cellA.Value = $row.date
    $cellA.NumberFormat = $dateFmt
    $ws.Cells.Item($row.r, 2).Value = $row.job
    $ws.Cells.Item($row.r, 3).Value = $row.cnt
    $ws.Cells.Item($row.r, 4).Value = $row.hrs
    $ws.Cells.Item($row.r, 5).Value = $row.emp
}

Write-Output "timesheet_daily_summary: wrote $($rows.Count) data rows (A2:E42)"
